# Weekly data refresh: insert a new price record at row 332 (most recent
# week's observation for "Orégano" at Mercado Mayorista Lo Valledor de
# Santiago) and push all the subsequent historical rows down by one, so the
# existing row 367 becomes row 368 and the sheet's used range grows from
# A1:R367 to A1:R368.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 332 - this shifts rows
# 332..367 down to 333..368, preserving all of their data/formatting.
$ws.Rows(332).Insert()

# Populate the newly inserted row 332 with the new weekly record.
$ws.Range("A332").Value = 6
$ws.Range("B332").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C332").Value = "Metropolitana"
$ws.Range("D332").Value = 45212
$ws.Range("E332").Value = 13
$ws.Range("F332").Value = 100112029
$ws.Range("G332").Value = "Orégano"
$ws.Range("H332").Value = "Sin especificar"
$ws.Range("I332").Value = "Primera"
$ws.Range("J332").Value = 32
$ws.Range("K332").Value = 16000
$ws.Range("L332").Value = 16000
$ws.Range("M332").Value = 16000
$ws.Range("N332").Value = "$/docena de atados"
$ws.Range("O332").Value = "Región Metropolitana"
$ws.Range("P332").Value = 5333
$ws.Range("Q332").Value = 3
$ws.Range("R332").Value = "Hortaliza"
